$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the launch time string (B2) to the new date value
# (leading apostrophe keeps the quote-prefix / text formatting that was already applied to this cell)
$ws.Range("B2").Value = "'18 Mar 2018 16:00:00.000'"

# Update the selection/active cell shown in the worksheet view
$ws.Range("B3").Select()
